# Add a new prediction-history sheet "Pred_2025-09-17_21-30-07" at the end of
# the workbook, mirroring the layout of the existing Pred_* sheets (header
# row with bold/bordered/centered style + one data row for the new match
# prediction: FC Bayern Munchen vs FC Augsburg).

$wb = $excel.ActiveWorkbook

# Use the last existing sheet as both the formatting template (for the
# header style) and the anchor to insert the new sheet after, so it lands
# at the end of the tab strip (matching the position in the diff).
$sheetCount = $wb.Worksheets.Count
$templateSheet = $wb.Worksheets.Item($sheetCount)
$templateHeaderRange = $templateSheet.Range("A1:R1")

$newSheet = $wb.Worksheets.Add($null, $templateSheet)
$newSheet.Name = "Pred_2025-09-17_21-30-07"

# Copy the header formatting (bold font, thin border, centered alignment)
# from the template sheet's header row so the new sheet reuses the same
# cell style rather than creating a near-duplicate style entry.
$templateHeaderRange.Copy()
$newSheet.Range("A1:R1").PasteSpecial(-4122)

$headers = @(
    "timestamp", "home_team", "away_team", "prediction", "prob_draw",
    "prob_homewin", "prob_awaywin", "home_rating", "away_rating",
    "home_form", "away_form", "home_momentum", "away_momentum",
    "is_derby", "home_value_eur", "away_value_eur", "home_avg_age",
    "away_avg_age"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$newSheet.Cells.Item(2, 1).Value = "2025-09-17T21:30:07"
$newSheet.Cells.Item(2, 2).Value = "fc bayern munchen"
$newSheet.Cells.Item(2, 3).Value = "fc augsburg"
$newSheet.Cells.Item(2, 4).Value = "HomeWin"
$newSheet.Cells.Item(2, 5).Value = 0.1645412378378072
$newSheet.Cells.Item(2, 6).Value = 0.7439006326607457
$newSheet.Cells.Item(2, 7).Value = 0.09155812950144707
$newSheet.Cells.Item(2, 8).Value = 31.13839876435202
$newSheet.Cells.Item(2, 9).Value = 41.42537269251397
$newSheet.Cells.Item(2, 10).Value = 1
$newSheet.Cells.Item(2, 11).Value = 0.2
$newSheet.Cells.Item(2, 12).Value = 18
$newSheet.Cells.Item(2, 13).Value = -5
$newSheet.Cells.Item(2, 14).Value = 0
$newSheet.Cells.Item(2, 15).Value = 905150000
$newSheet.Cells.Item(2, 16).Value = 134825000
$newSheet.Cells.Item(2, 17).Value = 26.04
$newSheet.Cells.Item(2, 18).Value = 24.82758620689655

Write-Output "Added sheet 'Pred_2025-09-17_21-30-07' with header + data row"
